# Update countries & provincias Spain
# Applies the periodic COVID-19 data refresh: updated case counts for several
# countries (causing a few of them to swap rank/order in the table) plus a
# refreshed "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Row 4: Estados Unidos -------------------------------------------------
$ws.Cells.Item(4, 2).Value = 391665   # Casos totales
$ws.Cells.Item(4, 3).Value = 24661    # Nuevos casos
$ws.Cells.Item(4, 4).Value = 21561    # Casos activos
$ws.Cells.Item(4, 5).Value = 357543   # Recuperados
$ws.Cells.Item(4, 7).Value = 1690     # Muertes hoy
$ws.Cells.Item(4, 8).Value = 12561    # Muertes

# --- Row 25: Noruega --------------------------------------------------------
$ws.Cells.Item(25, 2).Value = 5907
$ws.Cells.Item(25, 3).Value = 42
$ws.Cells.Item(25, 5).Value = 5786

# --- Rows 51-54: Colombia overtakes Sudafrica, Argentina e Islandia --------
# Colombia moves up into row 51 with fresh numbers; Sudafrica/Argentina/
# Islandia each drop one place, keeping the values they had before.
$ws.Cells.Item(51, 1).Value = "Colombia"
$ws.Cells.Item(51, 2).Value = 1780
$ws.Cells.Item(51, 3).Value = 201
$ws.Cells.Item(51, 4).Value = 100
$ws.Cells.Item(51, 5).Value = 1630
$ws.Cells.Item(51, 6).Value = 76
$ws.Cells.Item(51, 7).Value = 4
$ws.Cells.Item(51, 8).Value = 50

$ws.Cells.Item(52, 1).Value = "Sudafrica"
$ws.Cells.Item(52, 2).Value = 1749
$ws.Cells.Item(52, 3).Value = 63
$ws.Cells.Item(52, 4).Value = 95
$ws.Cells.Item(52, 5).Value = 1641
$ws.Cells.Item(52, 6).Value = 7
$ws.Cells.Item(52, 7).Value = 1
$ws.Cells.Item(52, 8).Value = 13

$ws.Cells.Item(53, 1).Value = "Argentina"
$ws.Cells.Item(53, 2).Value = 1628
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 338
$ws.Cells.Item(53, 5).Value = 1234
$ws.Cells.Item(53, 6).Value = 96
$ws.Cells.Item(53, 7).Value = 3
$ws.Cells.Item(53, 8).Value = 56

$ws.Cells.Item(54, 1).Value = "Islandia"
$ws.Cells.Item(54, 2).Value = 1586
$ws.Cells.Item(54, 3).Value = 24
$ws.Cells.Item(54, 4).Value = 559
$ws.Cells.Item(54, 5).Value = 1021
$ws.Cells.Item(54, 6).Value = 11
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 6

# --- Row 144: Islas Caimanes ------------------------------------------------
$ws.Cells.Item(144, 4).Value = 6
$ws.Cells.Item(144, 5).Value = 38

# --- Rows 182-183: Zimbabue overtakes Groenlandia --------------------------
# Zimbabue moves up into row 182 with fresh numbers; Groenlandia drops to 183
# keeping the values it had before.
$ws.Cells.Item(182, 1).Value = "Zimbabue"
$ws.Cells.Item(182, 2).Value = 11
$ws.Cells.Item(182, 3).Value = 1
$ws.Cells.Item(182, 4).Value = 0
$ws.Cells.Item(182, 5).Value = 9
$ws.Cells.Item(182, 6).Value = 0
$ws.Cells.Item(182, 7).Value = 1
$ws.Cells.Item(182, 8).Value = 2

$ws.Cells.Item(183, 1).Value = "Groenlandia"
$ws.Cells.Item(183, 2).Value = 11
$ws.Cells.Item(183, 3).Value = 0
$ws.Cells.Item(183, 4).Value = 10
$ws.Cells.Item(183, 5).Value = 1
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 0

# --- Title timestamp ---------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 7 de Abril de 2020 a las 22:52"
